$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.076.34'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.748.95'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5287'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2804'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06190'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("D10").Value = '1.747.21'
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07207'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.615'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9993'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9986'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '25.984.00'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006733'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("D21").Value = '1.970.91'
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.320'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.774'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.214'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.47%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08309'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.800'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.684'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04550'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.639'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("E35").Value = '  +4.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6331'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.93%  '
$ws.Range("E38").Value = '  +3.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.942'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9986'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3919'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7416'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.024'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1149'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.325'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05345'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.22%  '
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.683'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3465'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.64%  '
